$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume figures.
# Leading apostrophe forces Excel to treat the numeric-looking
# strings (e.g. "236.80") as literal text instead of silently
# re-parsing them as numbers (which would drop trailing zeros).
# Resetting the style back to "Normal" afterwards clears the
# quote-prefix formatting flag Excel applies, so the cell keeps
# its original (unstyled) appearance.

$ws.Range('D2').Value = "'42.914.43"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -1.43%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.341.81"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +1.22%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.01%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'307.07"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'100.78"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -1.50%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  -5.07%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D9').Value = "'0.512"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -3.66%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'34.94"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -2.56%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'52.15"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +0.24%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  -2.05%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = "'  -0.44%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = "'  -3.16%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'15.91"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +6.00%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'2.319.40"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +0.33%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'0.807"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -0.28%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'42.833.69"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -1.37%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'6.25"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +1.21%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  -1.98%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'11.72"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -5.86%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'67.95"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -0.38%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'236.80"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -2.20%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'2.02"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -0.51%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'2.56"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -2.48%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  -0.02%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'  +2.90%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  +1.00%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'35.18"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -4.45%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'9.35"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -3.24%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'160.10"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -4.17%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  -0.03%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  -3.17%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  +8.14%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'2.48"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -0.68%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'17.39"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -1.33%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'0.0728"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -2.41%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  -4.35%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  -0.42%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  -3.55%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'  -2.76%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'2.44"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +5.69%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'2.018.64"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +2.41%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = "'  -1.49%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'18.84"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -2.71%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'10.33"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +3.13%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'2.93"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -1.75%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'56.27"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +1.17%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'2.91"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -1.52%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'2.567.32"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +1.11%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = "'  +1.36%  "
$ws.Range('E51').Style = 'Normal'
